# Apply the two changes captured in the commit:
#   1. Slide 16's table switches from the custom "Table_0" style
#      ({E2E42672-6549-4812-BA5F-537ED8256C10}) to the built-in table
#      style {7D2FBF36-2037-4C1A-B9CA-98E2D3E22B1F}.
#   2. The presentation's theme colours revert from the "Integral" theme
#      back to the stock "Office Theme" palette (dk1/lt1 are identical in
#      both themes, so only the other ten slots need to change).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{7D2FBF36-2037-4C1A-B9CA-98E2D3E22B1F}")

# --- 2. Theme colours -------------------------------------------------------
# Office Theme RGB values packed the VBA way (R + G*256 + B*65536).
$officeThemeColors = @{
    3  = 6968388   # dk2       44546A
    4  = 15132391  # lt2       E7E6E6
    5  = 13998939  # accent1   5B9BD5
    6  = 3243501   # accent2   ED7D31
    7  = 10855845  # accent3   A5A5A5
    8  = 49407     # accent4   FFC000
    9  = 12874308  # accent5   4472C4
    10 = 4697456   # accent6   70AD47
    11 = 12673797  # hlink     0563C1
    12 = 7491477   # folHlink  954F72
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($slot in $officeThemeColors.Keys) {
    $themeColors.Item($slot).RGB = $officeThemeColors[$slot]
}
